# Revert "Powerpoint writer: consolidate text run nodes."
#
# The title text run "First " (on slide 1) and "Third " (on slide 3) each
# need to be split back into two separate runs: one for the word itself
# and one for the trailing space, so that the XML keeps:
#   <a:r><a:rPr/><a:t>First</a:t></a:r>
#   <a:r><a:rPr/><a:t> </a:t></a:r>
#   <a:r><a:rPr/><a:t>slide</a:t></a:r>
# instead of the single consolidated run <a:t>First </a:t>.

$p = $ppt.ActivePresentation

$slide1 = $p.Slides.Item(1)
$title1 = $slide1.Shapes.Item(1)
$tr1 = $title1.TextFrame.TextRange
# Locate the sub-range covering just "First" (not the trailing space that
# is still part of the same original run). Re-assigning this sub-range's
# own Text value forces the host to split the original run into two runs:
# "First" and " ", leaving the following "slide" run untouched.
$word1 = $tr1.Find("First", 0)
$word1.Text = "First"

$slide3 = $p.Slides.Item(3)
$title3 = $slide3.Shapes.Item(1)
$tr3 = $title3.TextFrame.TextRange
$word3 = $tr3.Find("Third", 0)
$word3.Text = "Third"
